# Updates sheet names and stimulus-file / condition values to reflect
# new run timestamps (improved accuracy of stimulus presentation time-logging)

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (new run timestamps baked into sheet names) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16512556020993557"
$wb.Worksheets.Item(2).Name = "NB_TO-16512556048793523"
$wb.Worksheets.Item(3).Name = "RS_TO-16512556048803525"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512556049273565"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512556050063517"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512556020623538.csv"
$ws1.Range("B3").Value = "GNG_stims-1651255602082353.csv"
$ws1.Range("B4").Value = "go_stims-16512556020833547.csv"
$ws1.Range("B5").Value = "GNG_stims-16512556020983543.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16512556037083514.csv"
$ws2.Range("B3").Value = "ZB-match_8-16512556025543523.csv"
$ws2.Range("B4").Value = "OB-16512556030483515.csv"
$ws2.Range("B5").Value = "ZB-match_1-16512556026583521.csv"
$ws2.Range("B6").Value = "TB-16512556046453533.csv"
$ws2.Range("B7").Value = "TB-1651255604851354.csv"
$ws2.Range("B8").Value = "OB-16512556038393536.csv"
$ws2.Range("B9").Value = "TB-16512556046193535.csv"
$ws2.Range("B10").Value = "ZB-match_6-16512556025793526.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512556048953514.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255604882352.csv"
$ws4.Range("B4").Value = "MM_stims-1651255604910354.csv"
$ws4.Range("B5").Value = "ZM_stims-16512556048953514.csv"
$ws4.Range("B6").Value = "MM_stims-1651255604926352.csv"
$ws4.Range("B7").Value = "ZM_stims-16512556049113548.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16512556049903557.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512556049743543.csv"
$ws5.Range("B4").Value = "SAT_stims-16512556049583514.csv"
$ws5.Range("B5").Value = "SAT_stims-16512556049323535.csv"
